# Update Name of Algo
# Apply updated RandomForest imputation values to Sheet1

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$updates = @{
    "D2"  = -6.897399999999996
    "C3"  = -10.85429999999998
    "A4"  = -21.61720000000003
    "B4"  = 4.947499999999999
    "C4"  = -10.83889999999999
    "E4"  = 13.38670000000002
    "B5"  = 4.983200000000001
    "A6"  = -21.26479999999999
    "A7"  = -21.53780000000001
    "B8"  = 4.651700000000004
    "C9"  = -11.74090000000001
    "C11" = -14.13950000000001
    "E12" = 11.14460000000001
    "C14" = -12.2154
    "A16" = -21.47600000000002
    "B16" = 4.996899999999997
    "E17" = 12.68540000000001
    "C18" = -14.58090000000001
    "A20" = -22.57690000000001
    "D20" = -8.166000000000007
    "E20" = 13.14009999999999
    "B22" = 5.499200000000001
    "C25" = -10.73299999999999
    "E25" = 13.46340000000001
}

foreach ($cellRef in $updates.Keys) {
    $ws.Range($cellRef).Value = $updates[$cellRef]
}
